$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("user_logs")

$rows = @(
    ,@(4, 3, 1, 'bW3KcumFHUEs4RaypuSXXQ==:uBFS5kqH6cHAIKnVr9HZKWjaf5c4aaBwzQgzeTMfnM9AOkGGEnbTeaNO3g3gfZ8ThrBPXX3zPS4JPr4QMQX/JTidRv7mqM0hACfcRmRgWmygxtqVFiOxdpCtZ1Q4mw8Uf2NWF3Jf4rC+Mq/hrwh+kA==', 45482.56212962963, 45482.56212962963)
    ,@(5, 4, 1, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.56240740741, 45482.56240740741)
    ,@(6, 5, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.57136574074, 45482.57136574074)
    ,@(7, 6, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.57188657407, 45482.57188657407)
    ,@(8, 7, 1, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.6, 45482.6)
    ,@(9, 8, 1, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.60024305555, 45482.60024305555)
    ,@(10, 9, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.92751157407, 45482.92751157407)
    ,@(11, 10, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.92756944444, 45482.92756944444)
    ,@(12, 11, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.92800925926, 45482.92800925926)
    ,@(13, 12, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.92826388889, 45482.92826388889)
    ,@(14, 13, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.92949074074, 45482.92949074074)
    ,@(15, 14, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.92952546296, 45482.92952546296)
    ,@(16, 15, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.93148148148, 45482.93148148148)
    ,@(17, 16, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.93155092592, 45482.93155092592)
    ,@(18, 17, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.93203703704, 45482.93203703704)
    ,@(19, 18, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.93207175926, 45482.93207175926)
    ,@(20, 19, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.93293981482, 45482.93293981482)
    ,@(21, 20, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.93313657407, 45482.93313657407)
    ,@(22, 21, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.93680555555, 45482.93680555555)
    ,@(23, 22, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.93702546296, 45482.93702546296)
    ,@(24, 23, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.93811342592, 45482.93811342592)
    ,@(25, 24, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.93824074074, 45482.93824074074)
    ,@(26, 25, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.93923611111, 45482.93923611111)
    ,@(27, 26, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.93930555556, 45482.93930555556)
    ,@(28, 27, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.94104166667, 45482.94104166667)
    ,@(29, 28, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.94114583333, 45482.94114583333)
    ,@(30, 29, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.94344907408, 45482.94344907408)
    ,@(31, 30, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.94362268518, 45482.94362268518)
    ,@(32, 31, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.9446875, 45482.9446875)
    ,@(33, 32, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.94583333333, 45482.94583333333)
    ,@(34, 33, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.94712962963, 45482.94712962963)
    ,@(35, 34, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.94723379629, 45482.94723379629)
    ,@(36, 35, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.94753472223, 45482.94753472223)
    ,@(37, 36, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.94770833333, 45482.94770833333)
    ,@(38, 37, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.94822916666, 45482.94822916666)
    ,@(39, 38, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.94829861111, 45482.94829861111)
    ,@(40, 39, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.94902777778, 45482.94902777778)
    ,@(41, 40, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.94914351852, 45482.94914351852)
    ,@(42, 41, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.94984953704, 45482.94984953704)
    ,@(43, 42, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.94994212963, 45482.94994212963)
    ,@(44, 43, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.95493055556, 45482.95493055556)
    ,@(45, 44, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.95545138889, 45482.95545138889)
    ,@(46, 45, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.95664351852, 45482.95664351852)
    ,@(47, 46, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.95706018519, 45482.95706018519)
    ,@(48, 47, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.95734953704, 45482.95734953704)
    ,@(49, 48, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.95780092593, 45482.95780092593)
    ,@(50, 49, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.97925925926, 45482.97925925926)
    ,@(51, 50, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.97934027778, 45482.97934027778)
    ,@(52, 51, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.98368055555, 45482.98368055555)
    ,@(53, 52, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.98399305555, 45482.98399305555)
    ,@(54, 53, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.9843287037, 45482.9843287037)
    ,@(55, 54, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.98539351852, 45482.98539351852)
    ,@(56, 55, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.98569444445, 45482.98569444445)
    ,@(57, 56, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.98583333333, 45482.98583333333)
    ,@(58, 57, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.98895833334, 45482.98895833334)
    ,@(59, 58, 2, 'bW3KcumFHUEs4RaypuSXXQ==:qCsQPaYdKSudUSJTF4cHsQ==', 45482.98909722222, 45482.98909722222)
    ,@(60, 59, 2, 'bW3KcumFHUEs4RaypuSXXQ==:c/twMdskw/rg9o8jlKaWIw==', 45482.99136574074, 45482.99136574074)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value2 = $row[5]
    $ws.Cells.Item($r, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Output "done"